# Update the "Divide and Conquer" sheet:
#  - add a "Comments" column (C)
#  - mark "Binary Search" as "Not Done" with a comment that the problem
#    is not available
#  - mark "Sum of Middle Elements of two sorted arrays" as "Done"
#    (solution was added)
#  - make this sheet the active tab / selected sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Divide and Conquer")

$ws.Range("C1").Value = "Comments"

$ws.Range("B4").Value = "Not Done"
$ws.Range("C4").Value = "Problem is not available"

$ws.Range("B5").Value = "Done"

$ws.Columns("C").ColumnWidth = 20

$ws.Activate()
$ws.Range("B6").Select()
